$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up rows 3-5 with the same base format as row 2 ---
$ws.Range("A2:D2").Copy($ws.Range("A3:D3"))
$ws.Range("A2:D2").Copy($ws.Range("A4:D4"))
$ws.Range("A2:D2").Copy($ws.Range("A5:D5"))

# Text values are entered in the order that reproduces the shared-string table
# order: row5's text first, then row3's, then row4's.
$ws.Range("A5").Value = "Definición del modelo de calidad, clasificación de issues"
$ws.Range("A3").Value = "Definición de estrategias de git branching y estándares de nomenclatura"
$ws.Range("A4").Value = "Definición del modelo de calidad y creación de issues"

# --- Row 3: Definición de estrategias de git branching y estándares de nomenclatura ---
$ws.Range("B3").Value = 45757
$ws.Range("C3").Value = 2
$ws.Rows.Item(3).RowHeight = 55.8

# --- Row 4: Definición del modelo de calidad y creación de issues ---
$ws.Range("B4").Value = 45759
$ws.Range("C4").Value = 3
$ws.Rows.Item(4).RowHeight = 42

# --- Row 5: Definición del modelo de calidad, clasificación de issues ---
$ws.Range("B5").Value = 45760
$ws.Range("C5").Value = 1.5
$ws.Rows.Item(5).RowHeight = 42

# --- Rows 6-8: blank, pre-formatted rows ---
$ws.Range("C2").Copy()
$ws.Range("B6:D8").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)

# --- Final selection matches the authored session's last click ---
$null = $ws.Range("C13").Select()
